# Generate Report for Handoff
#
# Updates the localization-status report after a new handoff generation run:
#   - Status text "Handed back: in sync with en-US" -> "Ready for handoff"
#   - Refreshed "Latest HO Xliff Generate Date" / handoff timestamps
#   - Priority "ht" -> "mt" for the a44599e7 file
#   - New Error Detail message recorded for the a44599e7 file (stale handback)
#   - A couple of report columns are narrowed / widened for readability

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Overview
$ws2 = $wb.Worksheets.Item(2)   # zh-cn
$ws3 = $wb.Worksheets.Item(3)   # de-de

$readyForHandoff = "Ready for handoff"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/19359c72d7d05967dd2186bb1a8bd89be9462f40/e2e/a44599e7-4040-4699-95bc-e3ab48235c91.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bc90dc50159feed881f0961d1a7dea995005bb7c/e2e/a44599e7-4040-4699-95bc-e3ab48235c91.md."

# ---- Overview sheet ----
$ws1.Range("E2").Value = $readyForHandoff
$ws1.Range("F2").Value = $readyForHandoff
$ws1.Range("G2").Value = "2016-11-02 05:17:00"
$ws1.Range("E3").Value = $readyForHandoff
$ws1.Range("F3").Value = $readyForHandoff
$ws1.Range("G3").Value = "2016-11-02 05:17:00"

# ---- zh-cn sheet ----
$ws2.Range("C2").Value = $readyForHandoff
$ws2.Range("E2").Value = "mt"
$ws2.Range("H2").Value = "2016-11-02 05:16:45"
$ws2.Range("P2").Value = $errorDetail
$ws2.Range("C3").Value = $readyForHandoff
$ws2.Range("E3").Value = "mt"
$ws2.Range("H3").Value = "2016-11-02 05:16:45"

# ---- de-de sheet ----
$ws3.Range("C2").Value = $readyForHandoff
$ws3.Range("E2").Value = "mt"
$ws3.Range("H2").Value = "2016-11-02 05:17:00"
$ws3.Range("P2").Value = $errorDetail
$ws3.Range("C3").Value = $readyForHandoff
$ws3.Range("E3").Value = "mt"
$ws3.Range("H3").Value = "2016-11-02 05:17:00"

# ---- column width tweaks ----
# Overview: zh-cn / de-de status columns narrower
$ws1.Columns.Item(5).ColumnWidth = 16.3333333333333
$ws1.Columns.Item(6).ColumnWidth = 16.3333333333333

# zh-cn / de-de: Status column narrower, Error Detail column wider
$ws2.Columns.Item(3).ColumnWidth = 16.3333333333333
$ws2.Columns.Item(16).ColumnWidth = 39.1666666666667

$ws3.Columns.Item(3).ColumnWidth = 16.3333333333333
$ws3.Columns.Item(16).ColumnWidth = 39.1666666666667
